# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (D3) and
# "Correspond Handback DateTime" (G3) timestamps on the per-language
# report sheets for the second (95e730fc...) file row.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D3").Value = "2016-02-24 12:04:18"
$zhcn.Range("G3").Value = "2016-02-24 12:05:02"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D3").Value = "2016-02-24 12:04:31"
$dede.Range("G3").Value = "2016-02-24 12:05:26"
